# feat: add 2022-Q3 data
#
# 1) Insert a brand-new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet (tab order becomes: 总计, 2022-Q3, 2022-Q2, 2022-Q1).
# 2) Populate it with the Q3 fund-holding table (same shape as the Q2 sheet).
# 3) Update the "总计" (summary) sheet: the row that used to describe 2022-Q2
#    now describes 2022-Q3 (with the new counts), a new row is inserted for
#    the 2022-Q2 summary that used to live there, and the 2022-Q1 row is
#    pushed down one row.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) New "2022-Q3" sheet, inserted before "2022-Q2"
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# match page margins used by the rest of the workbook
$wsQ3.PageSetup.LeftMargin = 54
$wsQ3.PageSetup.RightMargin = 54
$wsQ3.PageSetup.TopMargin = 72
$wsQ3.PageSetup.BottomMargin = 72
$wsQ3.PageSetup.HeaderMargin = 36
$wsQ3.PageSetup.FooterMargin = 36

# ---- header row (text, bold/centered style like the other sheets) ----
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$cols = @("B","C","D","E","F","G","H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsQ3.Range($cols[$i] + "1")
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
}
$wsQ2.Range("B1:H1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

# ---- data rows ----
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextCell $wsQ3.Range("B2") "159792"
Set-TextCell $wsQ3.Range("C2") "富国中证港股通互联网ETF"
Set-TextCell $wsQ3.Range("D2") "18.08"
Set-TextCell $wsQ3.Range("E2") "99.26"
Set-TextCell $wsQ3.Range("F2") "4.15"
Set-TextCell $wsQ3.Range("G2") "0.7503"
$wsQ3.Range("H2").Value = 6

Set-TextCell $wsQ3.Range("B3") "513770"
Set-TextCell $wsQ3.Range("C3") "华宝中证港股通互联网ETF"
Set-TextCell $wsQ3.Range("D3") "3.80"
Set-TextCell $wsQ3.Range("E3") "98.21"
Set-TextCell $wsQ3.Range("F3") "4.14"
Set-TextCell $wsQ3.Range("G3") "0.1573"
$wsQ3.Range("H3").Value = 6

Set-TextCell $wsQ3.Range("B4") "004321"
Set-TextCell $wsQ3.Range("C4") "前海开源沪港深强国产业灵活配置混合"
Set-TextCell $wsQ3.Range("D4") "0.11"
Set-TextCell $wsQ3.Range("E4") "78.52"
Set-TextCell $wsQ3.Range("F4") "5.43"
Set-TextCell $wsQ3.Range("G4") "0.0060"
$wsQ3.Range("H4").Value = 5

# column A (the little index numbers) carries the same style as the Q2 sheet
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("A4").Value = 2
$wsQ2.Range("A2:A4").Copy()
$wsQ3.Range("A2:A4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet
# ---------------------------------------------------------------------
# Row 2 used to describe 2022-Q2; it now describes 2022-Q3.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.91

# Row 3 is new: the 2022-Q2 summary that used to live in row 2.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 6
$wsTotal.Range("D3").Value = 0.52

# Row 4 is the old 2022-Q1 summary (previously row 3), shifted down.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0

$excel.CutCopyMode = 0
